$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New course "Joining Data with pandas" logged the same way as the other
# "Introduction/... " rows above it: copy row 17's look down into row 18,
# then overwrite with the new course name / rating.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "Joining Data with pandas"
$ws.Range("B18").Value = 3

$ws.Range("A19").Select()
